$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 579.6
$ws.Range("J9").Value = 575
$ws.Range("L9").Value = 575
$ws.Range("N9").Value = -913
$ws.Range("H19").Value = 5251
$ws.Range("J19").Value = 5787.357
$ws.Range("L19").Value = 5787.357
$ws.Range("N19").Value = -6137.357
$ws.Range("H76").Value = 4250.95
$ws.Range("I76").Value = 4056.611
$ws.Range("K76").Value = 4056.611
$ws.Range("M76").Value = -3741.611
$ws.Range("H79").Value = 4250.95
$ws.Range("I79").Value = 4056.611
$ws.Range("K79").Value = 4056.611
$ws.Range("M79").Value = -2964.611
$ws.Range("H80").Value = 4591.7856
$ws.Range("I80").Value = 482.6
$ws.Range("J80").Value = 6874.6665
$ws.Range("K80").Value = 1447.8
$ws.Range("L80").Value = 20623.9995
$ws.Range("M80").Value = -449.8000000000002
$ws.Range("N80").Value = -22619.9995
$ws.Range("H83").Value = 4591.7856
$ws.Range("I83").Value = 482.6
$ws.Range("J83").Value = 6874.6665
$ws.Range("K83").Value = 4343.400000000001
$ws.Range("L83").Value = 61871.9985
$ws.Range("M83").Value = 648.5999999999995
$ws.Range("N83").Value = -71855.9985
$ws.Range("H106").Value = 3650
$ws.Range("I106").Value = 3650
$ws.Range("K106").Value = 3650
$ws.Range("M106").Value = -3019
$ws.Range("H113").Value = 3813.1667
$ws.Range("I113").Value = 3157.6667
$ws.Range("J113").Value = 4468.6665
$ws.Range("K113").Value = 3157.6667
$ws.Range("L113").Value = 4468.6665
$ws.Range("M113").Value = 96.33329999999978
$ws.Range("N113").Value = -10976.6665
$ws.Range("H137").Value = 11122.743
$ws.Range("I137").Value = 18232.264
$ws.Range("J137").Value = 2680.1875
$ws.Range("K137").Value = 54696.792
$ws.Range("L137").Value = 8040.5625
$ws.Range("M137").Value = -52146.792
$ws.Range("N137").Value = -13140.5625
$ws.Range("H138").Value = 16627.117
$ws.Range("I138").Value = 1395.3585
$ws.Range("J138").Value = 70446
$ws.Range("K138").Value = 4186.0755
$ws.Range("L138").Value = 211338
$ws.Range("M138").Value = 953.9245000000001
$ws.Range("N138").Value = -221618
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 5005882.5
$ws.Range("I10").Value = 7504824
$ws.Range("J10").Value = 8000
$ws.Range("K10").Value = 7504824
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = -7504654
$ws.Range("N10").Value = -8340
$ws.Range("H74").Value = 436613.56
$ws.Range("I74").Value = 546281.4399999999
$ws.Range("K74").Value = 546281.4399999999
$ws.Range("M74").Value = -545407.4399999999
$ws.Range("H77").Value = 436613.56
$ws.Range("I77").Value = 546281.4399999999
$ws.Range("K77").Value = 2731407.2
$ws.Range("M77").Value = -2727039.2
$ws.Range("H88").Value = 6309.273
$ws.Range("J88").Value = 8915.286
$ws.Range("L88").Value = 8915.286
$ws.Range("N88").Value = -9727.286
$ws.Range("H91").Value = 6309.273
$ws.Range("J91").Value = 8915.286
$ws.Range("L91").Value = 8915.286
$ws.Range("N91").Value = -11723.286
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 49999.5
$ws.Range("J81").Value = 49999.5
$ws.Range("L81").Value = 49999.5
$ws.Range("N81").Value = -52121.5
$ws.Range("H84").Value = 49999.5
$ws.Range("J84").Value = 49999.5
$ws.Range("L84").Value = 149998.5
$ws.Range("N84").Value = -160606.5
$ws.Range("H138").Value = 105000.5
$ws.Range("J138").Value = 105000.5
$ws.Range("L138").Value = 105000.5
$ws.Range("N138").Value = -115280.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 75379.44500000001
$ws.Range("J9").Value = 75379.44500000001
$ws.Range("L9").Value = 75379.44500000001
$ws.Range("N9").Value = -75715.44500000001
$ws.Range("H22").Value = 1114.1428
$ws.Range("I22").Value = 466
$ws.Range("K22").Value = 466
$ws.Range("M22").Value = -116
$ws.Range("H94").Value = 2412.2856
$ws.Range("I94").Value = 2399.5
$ws.Range("J94").Value = 2417.4
$ws.Range("K94").Value = 2399.5
$ws.Range("L94").Value = 2417.4
$ws.Range("M94").Value = -1948.5
$ws.Range("N94").Value = -3319.4
$ws.Range("H99").Value = 10776.385
$ws.Range("I99").Value = 7597
$ws.Range("J99").Value = 12763.5
$ws.Range("K99").Value = 7597
$ws.Range("L99").Value = 12763.5
$ws.Range("M99").Value = -6099
$ws.Range("N99").Value = -15759.5
$ws.Range("H118").Value = 49500
$ws.Range("J118").Value = 49500
$ws.Range("L118").Value = 49500
$ws.Range("N118").Value = -52814
$ws.Range("H126").Value = 10776.385
$ws.Range("I126").Value = 7597
$ws.Range("J126").Value = 12763.5
$ws.Range("K126").Value = 22791
$ws.Range("L126").Value = 38290.5
$ws.Range("M126").Value = -20321
$ws.Range("N126").Value = -43230.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 351.54544
$ws.Range("I2").Value = 380.77777
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 2284.66662
$ws.Range("L2").Value = 1320
$ws.Range("M2").Value = -2171.66662
$ws.Range("N2").Value = -1546
$ws.Range("H8").Value = 201.08333
$ws.Range("I8").Value = 201.08333
$ws.Range("K8").Value = 603.24999
$ws.Range("M8").Value = -464.24999
$ws.Range("H131").Value = 148165.83
$ws.Range("J131").Value = 1882
$ws.Range("L131").Value = 5646
$ws.Range("N131").Value = -15726
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 400000320
$ws.Range("J14").Value = 533
$ws.Range("L14").Value = 533
$ws.Range("N14").Value = -869
$ws.Range("H26").Value = 11000
$ws.Range("I26").Value = 11000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -10720
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 11000
$ws.Range("I50").Value = 11000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 11000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -10502
$ws.Range("N50").ClearContents()
$ws.Range("H113").Value = 3772.5
$ws.Range("I113").Value = 2545
$ws.Range("K113").Value = 2545
$ws.Range("M113").Value = -375
$ws.Range("H126").Value = 3052.2856
$ws.Range("I126").Value = 2050.1667
$ws.Range("K126").Value = 6150.500100000001
$ws.Range("M126").Value = -3680.500100000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000.25
$ws.Range("I7").Value = 2998.6667
$ws.Range("K7").Value = 2998.6667
$ws.Range("M7").Value = -2886.6667
$ws.Range("H14").Value = 14004.5
$ws.Range("J14").Value = 14004.5
$ws.Range("L14").Value = 14004.5
$ws.Range("N14").Value = -14348.5
$ws.Range("H126").Value = 4000.25
$ws.Range("I126").Value = 2998.6667
$ws.Range("K126").Value = 8996.000100000001
$ws.Range("M126").Value = -6526.000100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 18500
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H132").Value = 1926.85
$ws.Range("I132").Value = 1317
$ws.Range("J132").Value = 3756.4
$ws.Range("K132").Value = 3951
$ws.Range("L132").Value = 11269.2
$ws.Range("M132").Value = -1421
$ws.Range("N132").Value = -16329.2
